$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the sheet with new blank rows 11-26 first, matching style of existing blank rows 9/10
# (do this before rows 9/10 get their new content so the style copied is still the "blank" style)
$ws.Range("A9:E10").Copy()
$ws.Range("A11:E26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in "Paid" for E5:E9 (copy formatting from an existing "Paid" cell, then set value)
$ws.Range("E2").Copy()
$ws.Range("E5:E9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("E5").Value = "Paid"
$ws.Range("E6").Value = "Paid"
$ws.Range("E7").Value = "Paid"
$ws.Range("E8").Value = "Paid"
$ws.Range("E9").Value = "Paid"

# A9 gets a new date (copy date formatting from A2, then set value)
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A9").Value = 45202

# Row 10: new transaction - date, formula amount, mode of payment
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A10").Value = 45204

$ws.Range("B10").Formula = "=900*5"

$ws.Range("C5").Copy()
$ws.Range("C10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C10").Value = "NEFT"

# Update selection to A16 as per the edit
$ws.Range("A16").Select() | Out-Null
